$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.693.36'

$ws.Range('D3').Value = '1.742.81'
$ws.Range('E3').Value = '  -5.46%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.12'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -8.64%  '

$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4924'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -6.67%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '41.58'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -7.59%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2435'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -23.04%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.05957'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -12.36%  '

$ws.Range('E11').Value = '  -5.34%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06777'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -12.96%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.70'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -23.11%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.466'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -10.93%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '76.99'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -12.81%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.5798'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -26.01%  '

$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.000'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.02%  '

$ws.Range('D19').Value = '25.736.94'
$ws.Range('E19').Value = '  -3.39%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.49'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -17.43%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.000006448'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -18.65%  '

$ws.Range('D22').Value = '1.966.69'
$ws.Range('E22').Value = '  -5.19%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.979'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -13.72%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.870'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -15.67%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.008'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -16.48%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '136.14'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -4.81%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.487'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -11.74%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.839'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -17.42%  '

$ws.Range('E29').Value = '  -14.66%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '100.78'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -9.11%  '

$ws.Range('E31').Value = '  -10.38%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08101'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -6.96%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.352'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -17.84%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04412'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -9.29%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9995'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.659'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -7.18%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.016'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -10.80%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.6081'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -16.82%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.714'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -12.58%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.060'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -12.11%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '103.56'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.27%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.01497'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -13.66%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7741'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -14.47%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.168'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -12.61%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.3752'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -22.29%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05119'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -12.14%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.1076'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -13.46%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.956'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -22.94%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '30.24'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -13.39%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '52.69'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -12.27%  '
